$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to Text format so that
# numeric-looking strings (e.g. "0.730", "3.70", thousand-dot prices)
# are preserved exactly as text instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "39.935.63"
$ws.Range("E2").Value = "  +0.31%  "

$ws.Range("D3").Value = "2.216.23"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "291.27"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").Value = "86.67"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("E7").Value = "  -0.42%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "0.467"
$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("D10").Value = "30.36"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "50.22"
$ws.Range("E11").Value = "  +5.77%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.0779"
$ws.Range("E12").Value = "  -0.82%  "

$ws.Range("E13").Value = "  +2.69%  "

$ws.Range("D14").Value = "6.42"
$ws.Range("E14").Value = "  +1.01%  "

$ws.Range("D15").Value = "2.560.39"
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").Value = "13.78"
$ws.Range("E16").Value = "  -2.06%  "

$ws.Range("D17").Value = "2.208.23"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").Value = "0.730"
$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("D19").Value = "39.870.89"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("D21").Value = "11.07"
$ws.Range("E21").Value = "  -3.65%  "

$ws.Range("D22").Value = "5.74"
$ws.Range("E22").Value = "  -1.27%  "

$ws.Range("D23").Value = "65.73"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "237.64"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  -0.71%  "

$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("D28").Value = "23.01"
$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("D29").Value = "9.22"
$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("E30").Value = "  -7.52%  "

$ws.Range("D31").Value = "156.22"
$ws.Range("E31").Value = "  +2.79%  "

$ws.Range("D32").Value = "31.88"
$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").Value = "4.96"
$ws.Range("E34").Value = "  +0.46%  "

$ws.Range("E35").Value = "  +6.25%  "

$ws.Range("D36").Value = "0.0713"
$ws.Range("E36").Value = "  -0.69%  "

$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -1.55%  "

$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("D39").Value = "0.0989"
$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("E40").Value = "  +2.45%  "

$ws.Range("D41").Value = "15.26"
$ws.Range("E41").Value = "  -4.38%  "

$ws.Range("D42").Value = "2.098.14"
$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").Value = "3.70"
$ws.Range("E43").Value = "  -2.15%  "

$ws.Range("D44").Value = "0.0271"
$ws.Range("E44").Value = "  +0.88%  "

$ws.Range("D45").Value = "17.86"
$ws.Range("E45").Value = "  +1.84%  "

$ws.Range("D46").Value = "9.78"
$ws.Range("E46").Value = "  -2.07%  "

$ws.Range("D47").Value = "2.00"
$ws.Range("E47").Value = "  -7.79%  "

$ws.Range("D48").Value = "2.71"
$ws.Range("E48").Value = "  +3.64%  "

$ws.Range("D49").Value = "2.434.26"
$ws.Range("E49").Value = "  +0.21%  "

$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("E51").Value = "  +2.51%  "

